$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("汽車")

# --- Header row (row 1): switch from "copy of row 2" to real column labels,
# and extend with the new property/legislator columns (matches the other
# sheets' schema: name/capacity/owner/register_date/register_reason/
# acquire_value/property_category/category/date/legislator_name/
# legislator_id/source_file/index).
$ws.Range("B1").Value2 = "name"
$ws.Range("C1").Value2 = "capacity"
$ws.Range("D1").Value2 = "owner"
$ws.Range("E1").Value2 = "register_date"
$ws.Range("F1").Value2 = "register_reason"
$ws.Range("G1").Value2 = "acquire_value"
$ws.Range("H1").Value2 = "property_category"
$ws.Range("I1").Value2 = "category"
$ws.Range("J1").Value2 = "date"
$ws.Range("K1").Value2 = "legislator_name"
$ws.Range("L1").Value2 = "legislator_id"
$ws.Range("M1").Value2 = "source_file"
$ws.Range("N1").Value2 = "index"

# New header cells get the same bold/centered/bordered look as the rest of
# row 1.
$hdrRange = $ws.Range("H1:N1")
$hdrRange.Font.Bold = $true
$hdrRange.HorizontalAlignment = -4108
$hdrRange.VerticalAlignment = -4160
$hdrRange.Borders.Item(7).LineStyle = 1
$hdrRange.Borders.Item(8).LineStyle = 1
$hdrRange.Borders.Item(9).LineStyle = 1
$hdrRange.Borders.Item(10).LineStyle = 1

# --- Row 2 (record #47, TOYOTACAMRYLE / 饒月琴) new trailing columns.
$ws.Range("H2").Value2 = "land"
$ws.Range("I2").Value2 = "normal"
$ws.Range("J2").Value2 = "2012-04-23"
$ws.Range("K2").Value2 = "許忠信"
$ws.Range("L2").Value2 = 1749
$ws.Range("M2").Value2 = "tmpa22c1"
$ws.Range("N2").Value2 = 47

# --- Row 3 (record #48, 中華GL20SS58 / 許忠信) new trailing columns.
$ws.Range("H3").Value2 = "land"
$ws.Range("I3").Value2 = "normal"
$ws.Range("J3").Value2 = "2012-04-23"
$ws.Range("K3").Value2 = "許忠信"
$ws.Range("L3").Value2 = 1749
$ws.Range("M3").Value2 = "tmpa22c1"
$ws.Range("N3").Value2 = 48
